# Applies the "Added Delete function in transactions.html" change:
# three transaction rows were removed from the "Spent" sheet (the web app's
# new Delete feature removed these rows from the exported database), the
# Txn ID column was auto-fit to its new content, and the last selected
# cell was moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three removed transactions, working from the bottom row
# upward so earlier deletions don't shift the row numbers of rows still
# to be removed.
#   Row 12: 20241205201916836536 / 2024-12-10 / 3333       / RDS PETROZA
#   Row 8:  20241205173730853658 / 2018-09-24 / 3023.25    / Berjaya Starbucks Coffee Company Sdn Bhd
#   Row 6:  20241205163526997030 / 2024-08-22 / 73.96      / RDS PETROZA (Others)
$ws.Range("A12").EntireRow.Delete()
$ws.Range("A8").EntireRow.Delete()
$ws.Range("A6").EntireRow.Delete()

# Resize column A (Txn ID) to fit its remaining content.
$ws.Columns.Item(1).AutoFit()

# Move the active selection, as last left by the editor.
$ws.Range("A24").Select()
